# The "codeforiati:group-code" and "codeforiati:group-name" columns (D and E)
# were swapped: the codelist generator now emits the group-code column before
# the group-name column, so for every data row the value that used to be in
# column D (group-name) now belongs in column E, and vice versa.
#
# Swap the contents of columns D and E for every used row on the sheet
# (headers included), which reproduces the effect of the upstream diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2

    $dCell.Value = $eVal
    $eCell.Value = $dVal
}
